$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the cell text uses a *literal* backslash-r-backslash-n marker (not an
# actual CR/LF), so it must be written with single-quoted (non-expanding)
# PowerShell strings.

# Row 30 (A30): small wording tweak - "우리가" -> "우리의"
$ws.Range("A30").Value = '우리는 아직\r\n우리의 가치를 몰라요\r\n당신은 특별해요!'

# Row 8 (A8): small wording tweak - remove the space between "고민" and "하고"
$ws.Range("A8").Value = '고민하고 있다면\r\n한번 해보는 게 어때요?'

# Update the active selection to A15 (matches the saved cursor position in the file)
$ws.Range("A15").Select()
